# Apply edits to "Uren log.xlsx":
# - Add a new row to the "Thomas" sheet logging "Searching data" hours on 2022-09-28.
# - Widen column A to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thomas")

# New entry: Wat?!? | Wanneer | Hoe lang
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Searching data"
$ws.Range("B4").Value = (Get-Date -Year 2022 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C4").Value = 1

# Widen column A (holds "Searching data") to fit its text, ~12.55 chars wide
$ws.Columns.Item(1).ColumnWidth = 11.6

$ws.Range("E7").Select()
